$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.048086285591125
$ws.Range("B1").Value = 1.73201310634613
$ws.Range("D1").Value = 2.696845293045044
$ws.Range("E1").Value = 1.232964396476746
